# "loss" becomes a new input column, inserted right before the existing
# column L (maxrgr01), pushing it and everything to its right one column
# over (L -> M, M -> N, ... AR -> AS). This mirrors Excel's "Insert" on
# the column, which shifts cells and carries formatting from the left
# neighbour (column K) into the new column.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at L; this shifts L:AR to M:AS and copies the
# formatting (style) of column K into the new column L.
$ws.Columns("L").Insert()

# Match column L's width to its neighbour K (both use the same style);
# the original template shows column L losing its "bestFit" flag but
# keeping the same numeric width (8) and style as K.
$ws.Columns("L").ColumnWidth = $ws.Columns("K").ColumnWidth()

# Header for the new column.
$ws.Cells.Item(1, 12).Value = "loss"

# Every data row (2-7) gets a loss value of 0.05.
$ws.Range("L2:L7").Value = 0.05

# Match the final selection recorded in the saved workbook.
[void]$ws.Range("K14").Select()
